$wb = $excel.ActiveWorkbook

# --- Sheet "sets": update home_points for match 1 / set 3 from 10 to 11 ---
$sets = $wb.Worksheets.Item("sets")
$sets.Range("D4").Value = 11

# --- Sheet "rallies": append a new rally row (row 74) ---
$rallies = $wb.Worksheets.Item("rallies")

$rallies.Cells.Item(74, 1).Value = 73          # A74 rally_id
$rallies.Cells.Item(74, 2).Value = 1           # B74 match_id
$rallies.Cells.Item(74, 3).Value = 3           # C74 set_number
$rallies.Cells.Item(74, 4).Value = 11          # D74 rally_no
$rallies.Cells.Item(74, 5).Value = "NOS"       # E74 side
$rallies.Cells.Item(74, 6).Value = ""          # F74 position
$rallies.Cells.Item(74, 7).Value = 2           # G74 player_number
$rallies.Cells.Item(74, 8).Value = "LINHA"     # H74 action
$rallies.Cells.Item(74, 9).Value = "PONTO"     # I74 result
$rallies.Cells.Item(74, 10).Value = "NOS"      # J74 who_scored
$rallies.Cells.Item(74, 11).Value = 11         # K74 score_home
$rallies.Cells.Item(74, 12).Value = 0          # L74 score_away
$rallies.Cells.Item(74, 13).Value = "1 2 l"    # M74 raw_text
$rallies.Cells.Item(74, 14).Value = "FRENTE"   # N74 position_zone
$rallies.Cells.Item(74, 15).Value = "FRENTE"   # O74 pos_fb
$rallies.Cells.Item(74, 16).Value = "FRENTE"   # P74 frente_fundo
